# Apply attendance value updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: G3, H3 -> 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: D4, E4 -> 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: D5, E5 -> 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: H6 -> 1
$ws.Range("H6").Value = 1

# Row 7: H7 -> 1
$ws.Range("H7").Value = 1

# Row 8: H8 -> 1
$ws.Range("H8").Value = 1

# Row 9: H9 -> 1
$ws.Range("H9").Value = 1

# Row 10: H10 -> 1
$ws.Range("H10").Value = 1

# Row 11: H11 -> 1
$ws.Range("H11").Value = 1

# Row 12: G12 -> 1
$ws.Range("G12").Value = 1

# Row 13: H13 -> 1
$ws.Range("H13").Value = 1

# Row 14: H14 -> 1
$ws.Range("H14").Value = 1

# Row 15: H15 -> 1
$ws.Range("H15").Value = 1

# Row 16: H16 -> 1
$ws.Range("H16").Value = 1

# Row 17: H17 -> 1
$ws.Range("H17").Value = 1

# Row 18: H18 -> 1
$ws.Range("H18").Value = 1
